$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.3
$ws.Range("H2").Value = 2.1
$ws.Range("J2").Value = 3.25
$ws.Range("K2").Value = 3.7
$ws.Range("L2").Value = 1.46
$ws.Range("N2").Value = 3.25
$ws.Range("O2").Value = 1.37
$ws.Range("P2").Value = 1.75
$ws.Range("Q2").Value = 2.12
$ws.Range("U2").Value = 1.94
$ws.Range("W2").Value = 1.32
$ws.Range("Y2").Value = 8.6
$ws.Range("AA2").Value = 29
$ws.Range("AF2").Value = 30
$ws.Range("AL2").Value = 70
$ws.Range("AN2").Value = 70
$ws.Range("F3").Value = 1.9
$ws.Range("H3").Value = 4.4
$ws.Range("I3").Value = 4.8
$ws.Range("K3").Value = 3.8
$ws.Range("L3").Value = 1.44
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 3.55
$ws.Range("P3").Value = 1.86
$ws.Range("Q3").Value = 2.04
$ws.Range("R3").Value = 1.32
$ws.Range("S3").Value = 3.75
$ws.Range("T3").Value = 1.86
$ws.Range("V3").Value = 1.27
$ws.Range("X3").Value = 15
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 980
$ws.Range("AA3").Value = 120
$ws.Range("AB3").Value = 9.199999999999999
$ws.Range("AN3").Value = 1000
$ws.Range("F4").Value = 2.68
$ws.Range("H4").Value = 2.78
$ws.Range("I4").Value = 3.1
$ws.Range("J4").Value = 3.05
$ws.Range("K4").Value = 3.6
$ws.Range("L4").Value = 1.48
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 3.2
$ws.Range("O4").Value = 1.39
$ws.Range("P4").Value = 1.72
$ws.Range("Q4").Value = 2.16
$ws.Range("R4").Value = 1.27
$ws.Range("S4").Value = 3.95
$ws.Range("V4").Value = 1.48
$ws.Range("X4").Value = 12
$ws.Range("Y4").Value = 10.5
$ws.Range("Z4").Value = 21
$ws.Range("AA4").Value = 55
$ws.Range("AB4").Value = 12.5
$ws.Range("AC4").Value = 7.8
$ws.Range("AD4").Value = 13.5
$ws.Range("AF4").Value = 22
$ws.Range("AG4").Value = 15.5
$ws.Range("AH4").Value = 22
$ws.Range("AJ4").Value = 50
$ws.Range("AN4").Value = 980
$ws.Range("AO4").Value = 44
$ws.Range("F5").Value = 1.51
$ws.Range("G5").Value = 1.58
$ws.Range("H5").Value = 7.2
$ws.Range("J5").Value = 4.2
$ws.Range("K5").Value = 4.9
$ws.Range("L5").Value = 1.4
$ws.Range("N5").Value = 3.85
$ws.Range("P5").Value = 1.97
$ws.Range("Q5").Value = 1.87
$ws.Range("R5").Value = 1.37
$ws.Range("X5").Value = 990
$ws.Range("Y5").Value = 29
$ws.Range("AC5").Value = 990
$ws.Range("AF5").Value = 9.4
$ws.Range("H6").Value = 2.36
$ws.Range("I6").Value = 2.46
$ws.Range("J6").Value = 3.5
$ws.Range("K6").Value = 3.7
$ws.Range("L6").Value = 1.39
$ws.Range("N6").Value = 3.8
$ws.Range("Q6").Value = 1.96
$ws.Range("S6").Value = 3.5
$ws.Range("T6").Value = 1.72
$ws.Range("V6").Value = 1.68
$ws.Range("AA6").Value = 40
$ws.Range("AF6").Value = 24
$ws.Range("AH6").Value = 17.5
$ws.Range("AO6").Value = 22
$ws.Range("H7").Value = 8.6
$ws.Range("J7").Value = 4.9
$ws.Range("L7").Value = 1.33
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 5.1
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 2.34
$ws.Range("Q7").Value = 1.69
$ws.Range("R7").Value = 1.53
$ws.Range("S7").Value = 2.78
$ws.Range("V7").Value = 1.11
$ws.Range("Z7").Value = 100
$ws.Range("AB7").Value = 12
$ws.Range("F8").Value = 1.67
$ws.Range("G8").Value = 1.76
$ws.Range("H8").Value = 6.6
$ws.Range("I8").Value = 8.4
$ws.Range("K8").Value = 3.75
$ws.Range("L8").Value = 1.58
$ws.Range("M8").Value = 1.12
$ws.Range("N8").Value = 2.62
$ws.Range("P8").Value = 1.52
$ws.Range("Q8").Value = 2.6
$ws.Range("R8").Value = 1.18
$ws.Range("S8").Value = 5.1
$ws.Range("T8").Value = 2.34
$ws.Range("U8").Value = 1.59
$ws.Range("W8").Value = 2.3
$ws.Range("G9").Value = 2.26
$ws.Range("H9").Value = 3.35
$ws.Range("I9").Value = 4.2
$ws.Range("J9").Value = 3.25
$ws.Range("L9").Value = 1.38
$ws.Range("N9").Value = 3.65
$ws.Range("O9").Value = 1.26
$ws.Range("P9").Value = 1.94
$ws.Range("Q9").Value = 1.81
$ws.Range("R9").Value = 1.36
$ws.Range("S9").Value = 3.05
$ws.Range("U9").Value = 2.12
$ws.Range("W9").Value = 1.79
